{"js": "// Split the \"Collision Detection system...\" paragraph into four paragraphs\n// with the follow-up thoughts, per the commit \"thoughts on next steps\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText =\n  \"- Collision Detection system- required to create a physics plugin to access low level API info\";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(targetText) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the Collision Detection paragraph\");\n}\n\n// Flat-OPC wrapped OOXML lets us write the exact run/paragraph structure\n// (separate <w:r> runs, <w:pPr><w:ind w:left=\"2160\"/></w:pPr>) instead of\n// having plain insertText calls silently coalesce into one run.\nconst flatOpc = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:tab/><w:t>- Collision Detection system</w:t></w:r>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:pPr><w:ind w:left=\"2160\"/></w:pPr>\n            <w:r><w:t>- obtain contact points and match to find which objects the collisions happened between</w:t></w:r>\n            <w:r><w:t xml:space=\"preserve\"> (completed)</w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:pPr><w:ind w:left=\"2160\"/></w:pPr>\n            <w:r><w:t>- score system providing rewards throughout the episode vs at the end of episode</w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:pPr><w:ind w:left=\"2160\"/></w:pPr>\n            <w:r><w:t xml:space=\"preserve\">- \\u2018debouncing\\u2019 so not penalised for contact which occurs over several time steps. </w:t></w:r>\n            <w:r><w:t>Using dictionary of switches to decide if this is a \\u2018new\\u2019 contact or not ( maybe this is not necessary to really discourage touching blocks of the wrong color?)</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntarget.getRange().insertOoxml(flatOpc, \"Replace\");\nawait context.sync();\n", "ps1": "# Split the \"Collision Detection system...\" paragraph into four paragraphs\n# with the follow-up thoughts, per the commit \"thoughts on next steps\".\n\n$d = $word.ActiveDocument\n\n$targetText = \"- Collision Detection system- required to create a physics plugin to access low level API info\"\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains($targetText)) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the Collision Detection paragraph\"\n}\n\n$wNs = \"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"\n\n$newXml = @\"\n<w:p xmlns:w='$wNs'><w:r><w:tab/><w:t>- Collision Detection system</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r></w:p><w:p xmlns:w='$wNs'><w:pPr><w:ind w:left=\"2160\"/></w:pPr><w:r><w:t>- obtain contact points and match to find which objects the collisions happened between</w:t></w:r><w:r><w:t xml:space=\"preserve\"> (completed)</w:t></w:r></w:p><w:p xmlns:w='$wNs'><w:pPr><w:ind w:left=\"2160\"/></w:pPr><w:r><w:t>- score system providing rewards throughout the episode vs at the end of episode</w:t></w:r></w:p><w:p xmlns:w='$wNs'><w:pPr><w:ind w:left=\"2160\"/></w:pPr><w:r><w:t xml:space=\"preserve\">- &#8216;debouncing&#8217; so not penalised for contact which occurs over several time steps. </w:t></w:r><w:r><w:t>Using dictionary of switches to decide if this is a &#8216;new&#8217; contact or not ( maybe this is not necessary to really discourage touching blocks of the wrong color?)</w:t></w:r></w:p>\n\"@\n\n$null = $target.Range.InsertXML($newXml)\n"}
